# Applies the commit "#327 Ajout des profils d'acces" changes:
#  1. Updates the "Date" metadata value (Metadata!B8) to the new timestamp.
#  2. On the "Elements" sheet, the two right-most "Mapping" columns (AK=37 and
#     AL=38) are swapped: the column that used to hold
#     "Mapping: RIM Mapping" / n/a / N/A / N/A values now holds
#     "Mapping: Spécification métier vers l'extension ROR ContactConfidentialityLevel"
#     / (empty) / (empty) / niveauConfidentialite, and vice versa. The column
#     widths of AK and AL are swapped as well.

$wb = $excel.ActiveWorkbook

# --- 1. Update Date on the Metadata sheet -------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- 2. Swap columns AK (37) and AL (38) on the Elements sheet ---------
$elements = $wb.Worksheets.Item("Elements")

# Row 1 (headers) - both cells are non-empty, safe straightforward swap.
$h1 = $elements.Cells.Item(1, 37).Value2
$h2 = $elements.Cells.Item(1, 38).Value2
$elements.Cells.Item(1, 37).Value2 = $h2
$elements.Cells.Item(1, 38).Value2 = $h1

# Row 2: both AK2 and AL2 already contain an empty value -> nothing to do.

# Row 3: AK3 = "n/a", AL3 = "" -> swap.
$elements.Cells.Item(3, 38).Value2 = "n/a"
$elements.Cells.Item(3, 37).Value2 = ""

# Row 4: both AK4 and AL4 already contain an empty value -> nothing to do.

# Row 5: AK5 = "N/A", AL5 = "" -> swap.
$elements.Cells.Item(5, 38).Value2 = "N/A"
$elements.Cells.Item(5, 37).Value2 = ""

# Row 6: AK6 = "N/A", AL6 = "niveauConfidentialite" -> swap.
$v6ak = $elements.Cells.Item(6, 37).Value2
$v6al = $elements.Cells.Item(6, 38).Value2
$elements.Cells.Item(6, 37).Value2 = $v6al
$elements.Cells.Item(6, 38).Value2 = $v6ak

# Swap the column widths too (target OOXML "width" = ColumnWidth + 5/6,
# so we compensate by subtracting 5/6 before assigning).
$offset = 5.0 / 6.0
$elements.Columns.Item(37).ColumnWidth = 83.7734375 - $offset
$elements.Columns.Item(38).ColumnWidth = 24.98046875 - $offset
